$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = "29.283.14"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "1.901.75"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("E4").Value = "  -0.53%  "
Set-TextValue $ws.Range("D5") "326.18"
$ws.Range("E5").Value = "  -0.65%  "
Set-TextValue $ws.Range("D6") "1.000"
$ws.Range("E6").Value = "  -0.35%  "
Set-TextValue $ws.Range("D7") "0.4650"
$ws.Range("E7").Value = "  +0.14%  "
Set-TextValue $ws.Range("D8") "0.3916"
$ws.Range("E8").Value = "  -0.36%  "
Set-TextValue $ws.Range("D9") "0.07884"
$ws.Range("E9").Value = "  -1.07%  "
Set-TextValue $ws.Range("D10") "0.9890"
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("E11").Value = "  -1.23%  "
$ws.Range("D12").Value = "1.932.87"
$ws.Range("E12").Value = "  +3.31%  "
Set-TextValue $ws.Range("D13") "7.081"
$ws.Range("E13").Value = "  -0.83%  "
Set-TextValue $ws.Range("D14") "5.753"
$ws.Range("E14").Value = "  -0.90%  "
Set-TextValue $ws.Range("D15") "0.06992"
$ws.Range("E15").Value = "  -0.06%  "
Set-TextValue $ws.Range("D16") "88.44"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("E17").Value = "  -0.31%  "
Set-TextValue $ws.Range("D19") "17.10"
$ws.Range("E19").Value = "  -1.00%  "
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").Value = "29.273.29"
$ws.Range("E21").Value = "  +0.27%  "
Set-TextValue $ws.Range("D22") "5.316"
$ws.Range("E22").Value = "  -1.23%  "
Set-TextValue $ws.Range("D23") "11.08"
$ws.Range("E23").Value = "  +0.03%  "
Set-TextValue $ws.Range("D24") "2.095"
$ws.Range("E24").Value = "  +1.78%  "
Set-TextValue $ws.Range("D25") "156.53"
$ws.Range("E25").Value = "  +0.03%  "
Set-TextValue $ws.Range("D26") "19.41"
$ws.Range("E26").Value = "  -0.94%  "
Set-TextValue $ws.Range("D27") "6.001"
$ws.Range("E27").Value = "  +2.51%  "
Set-TextValue $ws.Range("D28") "118.69"
$ws.Range("E28").Value = "  -0.94%  "
Set-TextValue $ws.Range("D29") "1.919"
$ws.Range("E29").Value = "  -4.46%  "
Set-TextValue $ws.Range("D30") "0.09377"
$ws.Range("E30").Value = "  -0.19%  "
Set-TextValue $ws.Range("D31") "0.9075"
$ws.Range("E31").Value = "  -1.91%  "
Set-TextValue $ws.Range("D32") "5.289"
$ws.Range("E32").Value = "  -1.62%  "
$ws.Range("E33").Value = "  -1.34%  "
Set-TextValue $ws.Range("D34") "3.212"
$ws.Range("E34").Value = "  -2.10%  "
Set-TextValue $ws.Range("D35") "1.185"
$ws.Range("E35").Value = "  +2.06%  "
Set-TextValue $ws.Range("D36") "0.05791"
$ws.Range("E36").Value = "  -0.97%  "
Set-TextValue $ws.Range("D37") "0.02091"
$ws.Range("E37").Value = "  -0.40%  "
$ws.Range("E38").Value = "  -0.39%  "
Set-TextValue $ws.Range("D39") "7.763"
$ws.Range("E39").Value = "  -3.11%  "
Set-TextValue $ws.Range("D40") "0.5719"
$ws.Range("E40").Value = "  -0.82%  "
Set-TextValue $ws.Range("D41") "0.1787"
$ws.Range("E41").Value = "  -1.36%  "
Set-TextValue $ws.Range("D42") "9.788"
$ws.Range("E42").Value = "  -2.22%  "
Set-TextValue $ws.Range("D43") "12.03"
$ws.Range("E43").Value = "  +0.05%  "
Set-TextValue $ws.Range("D44") "0.5350"
$ws.Range("E44").Value = "  -1.54%  "
Set-TextValue $ws.Range("D45") "2.194"
$ws.Range("E45").Value = "  -2.50%  "
Set-TextValue $ws.Range("D46") "0.07048"
$ws.Range("E46").Value = "  -0.91%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D47") "1.860"
$ws.Range("E47").Value = "  -1.45%  "
$ws.Range("B48").Value = "MXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D48") "2.590"
$ws.Range("E48").Value = "  +0.12%  "
Set-TextValue $ws.Range("D49") "113.36"
$ws.Range("E49").Value = "  +0.75%  "
$ws.Range("E50").Value = "  -1.98%  "
Set-TextValue $ws.Range("D51") "71.36"
$ws.Range("E51").Value = "  -0.62%  "
